# 229. Majority Element II
# Append a new tracking row to the LeetCode "Array" log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 46
$lastRow = $newRow - 1

# A46: problem number 229 - copy A45's formatting (centered style) first,
# then set the value so the pasted format sticks.
$ws.Cells.Item($lastRow, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item($newRow, 1).Value = 229

# B46: problem title (plain/unstyled cell, like the rest of the recent rows)
$ws.Cells.Item($newRow, 2).Value = "Majority Element II"
$ws.Cells.Item($newRow, 2).Style = "Normal"

# C46: language column
$ws.Cells.Item($newRow, 3).Value = "Java "

# D46: date solved (4/15/2023) - copy D45's date-number-format style too
$ws.Cells.Item($lastRow, 4).Copy()
$ws.Cells.Item($newRow, 4).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item($newRow, 4).Value = 45031

$ws.Application.CutCopyMode = $false

# Keep the worksheet's used range/dimension correct and restore the
# cursor/scroll position to roughly where the author left it.
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("E45").Select()

$wb.Save()
